$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- New issue row (row 44): add the new "Have Calibration set the mode
# properly" issue with status "Open" and status date 12/17/2018, matching
# the formatting already used by the surrounding rows.
$ws.Range("B44").Value = "Have Calibration set the mode properly"
$ws.Range("C44").Value = "Open"

# Copy the date formatting from the row above (D43) so the new date cell
# picks up the existing date style instead of minting a new one, then set
# the actual value.
$ws.Range("D43").Copy()
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("D44").Value = 43451

$excel.CutCopyMode = $false

# --- View state: scroll down toward the bottom of the list and move the
# selection to the newly added row.
$win = $excel.ActiveWindow
$win.ScrollRow = 32
$win.ScrollColumn = 1
$ws.Range("B45").Select()
